$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.44975366666667
$ws.Range("H2").Value = 46.349261
$ws.Range("I2").Value = 0.1287486886000874
$ws.Range("J2").Value = 0.1287486886000874
$ws.Range("M2").Value = 35.82261933333334
$ws.Range("N2").Value = 107.467858
$ws.Range("O2").Value = 0.1784748100644408
$ws.Range("P2").Value = 0.1784748100644408
$ws.Range("Q2").Value = 553.450644394771
$ws.Range("R2").Value = 4981.055799552938
$ws.Range("S2").Value = 0.02297839774394644
$ws.Range("T2").Value = 0.02297839774394644
$ws.Range("G3").Value = 15.44975366666667
$ws.Range("H3").Value = 46.349261
$ws.Range("I3").Value = 0.1287486886000874
$ws.Range("J3").Value = 0.1287486886000874
$ws.Range("O3").Value = 0.4121780644343741
$ws.Range("P3").Value = 0.4121780644343741
$ws.Range("Q3").Value = 1278.164774537246
$ws.Range("R3").Value = 11503.48297083521
$ws.Range("S3").Value = 0.053067385265648
$ws.Range("T3").Value = 0.05306738526564799
$ws.Range("G4").Value = 15.44975366666667
$ws.Range("H4").Value = 46.349261
$ws.Range("I4").Value = 0.1287486886000874
$ws.Range("J4").Value = 0.1287486886000874
$ws.Range("M4").Value = 73.74809799999998
$ws.Range("N4").Value = 221.244294
$ws.Range("O4").Value = 0.3674264480966141
$ws.Range("P4").Value = 0.3674264480966141
$ws.Range("Q4").Value = 1139.389947485192
$ws.Range("R4").Value = 10254.50952736673
$ws.Range("S4").Value = 0.04730567334942715
$ws.Range("T4").Value = 0.04730567334942715
$ws.Range("G5").Value = 15.44975366666667
$ws.Range("H5").Value = 46.349261
$ws.Range("I5").Value = 0.1287486886000874
$ws.Range("J5").Value = 0.1287486886000874
$ws.Range("M5").Value = 8.41412
$ws.Range("N5").Value = 25.24236
$ws.Range("O5").Value = 0.04192067740457094
$ws.Range("P5").Value = 0.04192067740457094
$ws.Range("Q5").Value = 129.9960813217733
$ws.Range("R5").Value = 1169.96473189596
$ws.Range("S5").Value = 0.005397232241065825
$ws.Range("T5").Value = 0.005397232241065824
$ws.Range("I6").Value = 0.4074352211478151
$ws.Range("J6").Value = 0.4074352211478151
$ws.Range("M6").Value = 35.82261933333334
$ws.Range("N6").Value = 107.467858
$ws.Range("O6").Value = 0.1784748100644408
$ws.Range("P6").Value = 0.1784748100644408
$ws.Range("Q6").Value = 1751.43753420127
$ws.Range("R6").Value = 15762.93780781143
$ws.Range("S6").Value = 0.07271692370791974
$ws.Range("T6").Value = 0.07271692370791974
$ws.Range("I7").Value = 0.4074352211478151
$ws.Range("J7").Value = 0.4074352211478151
$ws.Range("O7").Value = 0.4121780644343741
$ws.Range("P7").Value = 0.4121780644343741
$ws.Range("S7").Value = 0.1679358608350976
$ws.Range("T7").Value = 0.1679358608350976
$ws.Range("I8").Value = 0.4074352211478151
$ws.Range("J8").Value = 0.4074352211478151
$ws.Range("M8").Value = 73.74809799999998
$ws.Range("N8").Value = 221.244294
$ws.Range("O8").Value = 0.3674264480966141
$ws.Range("P8").Value = 0.3674264480966141
$ws.Range("Q8").Value = 3605.687951270609
$ws.Range("R8").Value = 32451.19156143549
$ws.Range("S8").Value = 0.1497024761358002
$ws.Range("T8").Value = 0.1497024761358002
$ws.Range("I9").Value = 0.4074352211478151
$ws.Range("J9").Value = 0.4074352211478151
$ws.Range("M9").Value = 8.41412
$ws.Range("N9").Value = 25.24236
$ws.Range("O9").Value = 0.04192067740457094
$ws.Range("P9").Value = 0.04192067740457094
$ws.Range("Q9").Value = 411.3826922634
$ws.Range("R9").Value = 3702.444230370601
$ws.Range("S9").Value = 0.01707996046899758
$ws.Range("T9").Value = 0.01707996046899758
$ws.Range("G10").Value = 24.32144666666666
$ws.Range("H10").Value = 72.96433999999999
$ws.Range("I10").Value = 0.2026798893205849
$ws.Range("J10").Value = 0.2026798893205849
$ws.Range("M10").Value = 35.82261933333334
$ws.Range("N10").Value = 107.467858
$ws.Range("O10").Value = 0.1784748100644408
$ws.Range("P10").Value = 0.1784748100644408
$ws.Range("Q10").Value = 871.2579255759688
$ws.Range("R10").Value = 7841.32133018372
$ws.Range("S10").Value = 0.03617325475037327
$ws.Range("T10").Value = 0.03617325475037327
$ws.Range("G11").Value = 24.32144666666666
$ws.Range("H11").Value = 72.96433999999999
$ws.Range("I11").Value = 0.2026798893205849
$ws.Range("J11").Value = 0.2026798893205849
$ws.Range("O11").Value = 0.4121780644343741
$ws.Range("P11").Value = 0.4121780644343741
$ws.Range("Q11").Value = 2012.123757169697
$ws.Range("R11").Value = 18109.11381452728
$ws.Range("S11").Value = 0.08354020447993184
$ws.Range("T11").Value = 0.08354020447993184
$ws.Range("G12").Value = 24.32144666666666
$ws.Range("H12").Value = 72.96433999999999
$ws.Range("I12").Value = 0.2026798893205849
$ws.Range("J12").Value = 0.2026798893205849
$ws.Range("M12").Value = 73.74809799999998
$ws.Range("N12").Value = 221.244294
$ws.Range("O12").Value = 0.3674264480966141
$ws.Range("P12").Value = 0.3674264480966141
$ws.Range("Q12").Value = 1793.660432275106
$ws.Range("R12").Value = 16142.94389047596
$ws.Range("S12").Value = 0.07446995183367737
$ws.Range("T12").Value = 0.07446995183367737
$ws.Range("G13").Value = 24.32144666666666
$ws.Range("H13").Value = 72.96433999999999
$ws.Range("I13").Value = 0.2026798893205849
$ws.Range("J13").Value = 0.2026798893205849
$ws.Range("M13").Value = 8.41412
$ws.Range("N13").Value = 25.24236
$ws.Range("O13").Value = 0.04192067740457094
$ws.Range("P13").Value = 0.04192067740457094
$ws.Range("Q13").Value = 204.6435708269333
$ws.Range("R13").Value = 1841.7921374424
$ws.Range("S13").Value = 0.008496478256602382
$ws.Range("T13").Value = 0.008496478256602382
$ws.Range("G14").Value = 31.33616366666666
$ws.Range("H14").Value = 94.00849099999999
$ws.Range("I14").Value = 0.2611362009315126
$ws.Range("J14").Value = 0.2611362009315126
$ws.Range("M14").Value = 35.82261933333334
$ws.Range("N14").Value = 107.467858
$ws.Range("O14").Value = 0.1784748100644408
$ws.Range("P14").Value = 0.1784748100644408
$ws.Range("Q14").Value = 1122.543462398031
$ws.Range("R14").Value = 10102.89116158228
$ws.Range("S14").Value = 0.04660623386220136
$ws.Range("T14").Value = 0.04660623386220136
$ws.Range("G15").Value = 31.33616366666666
$ws.Range("H15").Value = 94.00849099999999
$ws.Range("I15").Value = 0.2611362009315126
$ws.Range("J15").Value = 0.2611362009315126
$ws.Range("O15").Value = 0.4121780644343741
$ws.Range("P15").Value = 0.4121780644343741
$ws.Range("Q15").Value = 2592.454315584485
$ws.Range("R15").Value = 23332.08884026037
$ws.Range("S15").Value = 0.1076346138536967
$ws.Range("T15").Value = 0.1076346138536967
$ws.Range("G16").Value = 31.33616366666666
$ws.Range("H16").Value = 94.00849099999999
$ws.Range("I16").Value = 0.2611362009315126
$ws.Range("J16").Value = 0.2611362009315126
$ws.Range("M16").Value = 73.74809799999998
$ws.Range("N16").Value = 221.244294
$ws.Range("O16").Value = 0.3674264480966141
$ws.Range("P16").Value = 0.3674264480966141
$ws.Range("Q16").Value = 2310.982469033372
$ws.Range("R16").Value = 20798.84222130035
$ws.Range("S16").Value = 0.0959483467777094
$ws.Range("T16").Value = 0.0959483467777094
$ws.Range("G17").Value = 31.33616366666666
$ws.Range("H17").Value = 94.00849099999999
$ws.Range("I17").Value = 0.2611362009315126
$ws.Range("J17").Value = 0.2611362009315126
$ws.Range("M17").Value = 8.41412
$ws.Range("N17").Value = 25.24236
$ws.Range("O17").Value = 0.04192067740457094
$ws.Range("P17").Value = 0.04192067740457094
$ws.Range("Q17").Value = 263.6662414309733
$ws.Range("R17").Value = 2372.99617287876
$ws.Range("S17").Value = 0.01094700643790516
$ws.Range("T17").Value = 0.01094700643790516
